# This script reproduces the changes made to NIT-9017737532.xlsx:
#  - the second worker/period row in the detail table (old row 17, the
#    duplicate "73152213 / GUSTAVO ADOLFO LLERENA AVILA" line for period
#    2504) is removed entirely - "parte 1" of the new statement only has
#    one period line left, which shifts the signature rows up;
#  - the summary figures are refreshed for the new account statement:
#      "Cant. Periodos" (F13) goes from 2 to 1
#      "VALOR MORA" (E11) goes from 74022 to 17082
#      the remaining detail row's "Salario Basico" (G16) goes from
#      2611297 to 1423500

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second detail row (period 2504). Deleting the whole row
# shifts everything below it up by one and drops the now-unused
# "2504" shared string automatically.
$ws.Rows("17:17").Delete()

# Refresh the summary / detail values for the new statement.
$ws.Range("E11").Value = 17082
$ws.Range("F13").Value = 1
$ws.Range("G16").Value = 1423500
